$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price column (D) cells whose new values could
# otherwise be auto-parsed as numbers by Excel (losing formatting like
# trailing zeros, e.g. "12.60" -> 12.6). The Volume(1h) column (E) values
# already contain spaces/% so they stay text naturally.

$ws.Range('D2').Value = '22.999.20'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.577.62'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '299.54'
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3749'
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3557'
$ws.Range('E8').Value = '  -2.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '49.93'
$ws.Range('E9').Value = '  +2.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.003'
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.216'
$ws.Range('E11').Value = '  -3.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07972'
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.86'
$ws.Range('E13').Value = '  -4.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.402'
$ws.Range('E14').Value = '  -2.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.305'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001223'
$ws.Range('E16').Value = '  -3.02%  '
$ws.Range('D17').Value = '1.579.34'
$ws.Range('E17').Value = '  -1.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.16'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06736'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.65'
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.345'
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('D23').Value = '23.001.32'
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.60'
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.369'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.830'
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.58'
$ws.Range('E27').Value = '  -2.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.35'
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.162'
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.11'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.330'
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.538'
$ws.Range('E32').Value = '  -3.67%  '
$ws.Range('D33').Value = '1.753.20'
$ws.Range('E33').Value = '  -1.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9314'
$ws.Range('E34').Value = '  -3.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07328'
$ws.Range('E35').Value = '  -4.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08751'
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.935'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02631'
$ws.Range('E38').Value = '  -4.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2459'
$ws.Range('E39').Value = '  -2.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.979'
$ws.Range('E40').Value = '  -3.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.343'
$ws.Range('E41').Value = '  -3.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6862'
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.91'
$ws.Range('E43').Value = '  -6.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.67'
$ws.Range('E44').Value = '  -7.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6335'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.966'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.241'
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '129.97'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07850'
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.180'
$ws.Range('E51').Value = '  +1.16%  '
